# Penalty Reward System (unfinished) - shift Week_Start_Date by one week and
# update MyForecast values, plus refresh the derived Summary sheet metrics.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Forecast Comparison" ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# New Week_Start_Date values (col B) - each row's date moves one week later.
$weekStartDates = @{
    2  = "2025-01-12"
    3  = "2025-01-19"
    4  = "2025-01-26"
    5  = "2025-02-02"
    6  = "2025-02-09"
    7  = "2025-02-16"
    8  = "2025-02-23"
    9  = "2025-03-02"
    10 = "2025-03-09"
    11 = "2025-03-16"
    12 = "2025-03-23"
    13 = "2025-03-30"
    14 = "2025-04-06"
    15 = "2025-04-13"
    16 = "2025-04-20"
    17 = "2025-04-27"
}

# New MyForecast values (col D)
$myForecast = @{
    2  = 8
    3  = 9
    4  = 10
    5  = 12
    6  = 12
    7  = 12
    8  = 13
    9  = 13
    10 = 13
    11 = 13
    12 = 14
    13 = 14
    14 = 14
    15 = 13
    16 = 14
    17 = 13
}

foreach ($row in 2..17) {
    # Leading apostrophe forces the date-looking text to stay plain text,
    # matching the original inline-string storage instead of becoming a
    # real Excel date serial number.
    $ws1.Cells.Item($row, 2).Value = "'" + $weekStartDates[$row]
    $ws1.Cells.Item($row, 4).Value = $myForecast[$row]
}

# --- Sheet 2: "Summary" ---
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B2").Value  = "2023-02-19 to 2025-01-05"
$ws2.Range("B4").Value  = "'242"
$ws2.Range("B5").Value  = "'60"
$ws2.Range("B6").Value  = "'26"
$ws2.Range("B8").Value  = "3627 units"
$ws2.Range("B9").Value  = "'197"
$ws2.Range("B10").Value = "'90"
$ws2.Range("B11").Value = "'39"
$ws2.Range("B12").Value = "'14"
$ws2.Range("B13").Value = "'2025-03-23"
$ws2.Range("B14").Value = "'8"
$ws2.Range("B15").Value = "'2025-01-12"
